$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# New "fase2" indicators appended to the Library_Formula sheet, continuing
# the existing CREATE/MODIFY / LIB_EWS_IT / <indicator> / <blank> / String / String
# pattern used by the rows above (e.g. row 90 = INDICATOR_59).
$indicators = @(
    "INDICATOR_62",
    "INDICATOR_63",
    "INDICATOR_64",
    "INDICATOR_65",
    "INDICATOR_66",
    "INDICATOR_67",
    "INDICATOR_68",
    "INDICATOR_69",
    "INDICATOR_70"
)

$row = 91
foreach ($ind in $indicators) {
    foreach ($col in 1,2,3,5,6) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Font.Name = "Trebuchet MS"
        $cell.Font.Size = 10
    }

    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($row, 3).Value = $ind
    $ws.Cells.Item($row, 5).Value = "String"
    $ws.Cells.Item($row, 6).Value = "String"

    $row = $row + 1
}

# Reflect the scroll position / selection change recorded for this sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 75
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E90:F99").Select()
